$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# 1. Update Status text on Overview sheet (E2 and F2)
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# 2. Update zh-cn sheet: K2 (Latest Handback DateTime), P2 (Error Detail)
$zhcn.Range("K2").Value = "2016-09-04 16:54:29"
$zhcn.Range("P2").Value = ""

# 3. Update de-de sheet: K2 (Latest Handback DateTime), P2 (Error Detail)
$dede.Range("K2").Value = "2016-09-04 16:54:36"
$dede.Range("P2").Value = ""

# 4. Column width changes
$overview.Range("E1").ColumnWidth = 29.9777047293527
$overview.Range("F1").ColumnWidth = 29.9777047293527

$zhcn.Range("C1").ColumnWidth = 29.9777047293527
$zhcn.Range("P1").ColumnWidth = 13.7470528738839

$dede.Range("C1").ColumnWidth = 29.9777047293527
$dede.Range("P1").ColumnWidth = 13.7470528738839
